$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: log the actual pages written that day ("webmapper" progress in the
# design section) and flag the entry with the "(Design) 2.2 data viz" note.
$ws.Range("E14").Value = 44
$ws.Range("H14").Value = "(Design) 2.2 data viz"

# Row 15: same for the next day.
$ws.Range("E15").Value = 46
$ws.Range("H15").Value = "(Design) 2.2 data viz"

# Portrait page setup for printing.
$ws.PageSetup.Orientation = 1

# Leave the selection where it was left after the edit.
$ws.Range("H15").Select()
